# Remove the straight connector shape "Straight Connector 111" (id 112)
# from slide 4 of the presentation, per the Phase 1 workflow readme
# figure update.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Straight Connector 111") {
        $shape.Delete()
        break
    }
}
